# Update the extracted requisition record (row 2) with the new patient/record data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be treated as literal text so numeric/date-looking
    # strings (ids, zip codes, phone numbers, dates) keep their exact
    # formatting (leading zeros, dashes, etc.) instead of being coerced
    # into a number/date by Excel's input parser.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Core record fields
Set-TextValue "C2" "20474"
Set-TextValue "D2" "21843299"
Set-TextValue "E2" "MM0000012163"
Set-TextValue "F2" "Sagis DX"
Set-TextValue "G2" "Cardarelle Ross"
Set-TextValue "H2" "Joy"

# MiddleName no longer present for this patient
$ws.Range("I2").ClearContents()

Set-TextValue "J2" "1950-06-08"
Set-TextValue "K2" "Female"
Set-TextValue "L2" "5601 Spruce Ave"
Set-TextValue "M2" "CO"
Set-TextValue "N2" "Castle Rock"
Set-TextValue "O2" "80104"
Set-TextValue "P2" "3038600948"

Set-TextValue "R2" "Taylor Todd,"
Set-TextValue "S2" "Clarity Dermatology - Castle Rock"

# Primary insurance
Set-TextValue "T2" "Joy Cardarelle Ross"
Set-TextValue "U2" "Self"

# PrimaryInsurance_SubDOB no longer populated
$ws.Range("V2").ClearContents()

Set-TextValue "W2" "Medicare of Colorado"
Set-TextValue "Y2" "9VD5TU8PW68"

# Secondary insurance
Set-TextValue "AD2" "Joy Cardarelle Ross"
Set-TextValue "AE2" "Self"

# SecondaryInsurance_SubDOB no longer populated
$ws.Range("AF2").ClearContents()

Set-TextValue "AG2" "AARP Medicare Supplement/Fixed Indemnity by UHC"
Set-TextValue "AI2" "01477099411"
